$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the two new columns (L's text is entered first so it lands
# earlier in the shared-strings table, matching the authored workbook)
$ws.Range("L1").Value = "diff"
$ws.Range("K1").Value = "% growth"

# Per-row formulas: K = J/I (% growth), L = J-I (diff).
# Row 2 is entered on its own, then K3:K29 / L3:L29 are entered as one
# range-fill so the engine records them as a shared formula group, mirroring
# how the original author filled the formula down the column.
$ws.Range("K2").Formula = "=J2/I2"
$ws.Range("L2").Formula = "=J2-I2"
$ws.Range("K3:K29").Formula = "=J3/I3"
$ws.Range("L3:L29").Formula = "=J3-I3"

# Summary formulas
$ws.Range("M2").Formula = "=AVERAGE(L2:L29)"
$ws.Range("M3").Formula = "=MEDIAN(L2:L29)"

# Match new column width for K to the same as column J
$ws.Columns.Item(11).ColumnWidth = $ws.Columns.Item(10).ColumnWidth

# Update selection to match final state
$ws.Range("N10").Select()
